$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text/model (swap "hasName"/"hasGeoname" style headers
#     for plain "ID"/"Geoname"/"Name"/"Wonderland Location") ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Geoname"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Wonderland Location"

# --- Bump the whole sheet's base font size 10/12 -> 14 (theme color) ---
$ws.Cells.Font.Size = 14

# --- Remove the special larger-font styling that used to single out D8:D10 ---
$ws.Range("D8:D10").Font.Size = 14
$ws.Range("D8:D10").Font.ThemeColor = 1

# --- Re-select cursor like the author left it ---
$ws.Range("B8").Select
